# The "optimization_parameters" sheet had a stray row (row 16: a leftover
# "Sheet" label with values 3 / 4) sitting between the Strain list (row 15)
# and the simulation_timepoints row (row 17). Clean it up by selecting the
# whole row and deleting it, which shifts simulation_timepoints up to row 16
# and drops the now-unused "Sheet" shared string / number-format style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

$ws.Rows.Item(16).Select() | Out-Null
$ws.Rows.Item(16).Delete()
